# correcting errors in Zslgf
$wb = $excel.ActiveWorkbook

# --- Sheet "RVC": update measured line-length / load values (D4, D5, D6) and DP (B9) ---
$rvc = $wb.Worksheets.Item("RVC")
$rvc.Range("D4").Value = 6600
$rvc.Range("D5").Value = 3960
$rvc.Range("D6").Value = 1320
$rvc.Range("B9").Value = 922000

# --- Sheet "FaultsPOC": fix the weighted-average Zslgf formulas in row 5 ---
$faults = $wb.Worksheets.Item("FaultsPOC")
$faults.Range("B5").Formula = "=(2*B3+B4)/3"
$faults.Range("C5").Formula = "=(2*C3+C4)/3"

# --- Restore selections on each sheet ---
$faults.Range("C5").Select()

$lineCodes = $wb.Worksheets.Item("LineCodes")
$lineCodes.Range("B19").Select()

$rvc.Activate()
$rvc.Range("B10").Select()
